# completed the test_triangle.py file
# Fill in the developer name and the previously-blank "Preconditions",
# "Method Inputs" and "Expected Result" columns of the __init__ / __str__ /
# calculate_area / calculate_perimeter test rows of the unit-test plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name (row 3, next to the "Developer:" label)
$ws.Range("C3").Value = "Bibekdeep Singh"

# Test 1 - __init__ : Attribute set to input values.
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'color:"black",                                                                  length: 3,                                                                         width: 4'

# Test 2 - __init__ : Exception raised when color is blank
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'color:"",                                                                  length: 3,                                                                         width: 4'

$ws.Range("G7").Value = "attributes got set"
$ws.Range("G8").Value = "Value Error"

# Test 3 - __init__ : Exception raised when length is not an integer.
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'color:"black",                                                                  length: "three",                                                                         width: 4'
$ws.Range("G9").Value = "Value Error"

# Test 4 - __init__ : Exception raised when width is not an integer.
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'color:"black",                                                                  length: 3,                                                                         width: "four"'
$ws.Range("G10").Value = "Value Error"

# Test 5 - __str__ : Returns string formatted appropriately
$ws.Range("F11").Value = "None"

# Test 6 - calculate_area : Returns correct calculated value.
$ws.Range("F12").Value = "None"

# Test 7 - calculate_perimeter : Returns correct calculated value.
$ws.Range("F13").Value = "None"

# Leave the cursor where the author finished typing.
$ws.Range("E11").Select()
